$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 925
$ws.Range("I12").Value = 893.3333
$ws.Range("J12").Value = 972.5
$ws.Range("K12").Value = 893.3333
$ws.Range("L12").Value = 972.5
$ws.Range("M12").Value = -723.3333
$ws.Range("N12").Value = -1312.5
$ws.Range("H13").Value = 19999
$ws.Range("J13").Value = 19999
$ws.Range("L13").Value = 19999
$ws.Range("N13").Value = -20337
$ws.Range("H43").Value = 4143.3335
$ws.Range("J43").Value = 4756.2856
$ws.Range("L43").Value = 4756.2856
$ws.Range("N43").Value = -4894.2856
$ws.Range("H86").Value = 5642.857
$ws.Range("H89").Value = 5642.857
$ws.Range("H98").Value = 652.5
$ws.Range("I98").Value = 652.5
$ws.Range("K98").Value = 652.5
$ws.Range("M98").Value = 845.5
$ws.Range("H106").Value = 37624.25
$ws.Range("J106").Value = 24500
$ws.Range("L106").Value = 24500
$ws.Range("N106").Value = -25762
$ws.Range("H113").Value = 3665
$ws.Range("I113").Value = 3995
$ws.Range("K113").Value = 3995
$ws.Range("M113").Value = -741
$ws.Range("H122").Value = 652.5
$ws.Range("I122").Value = 652.5
$ws.Range("K122").Value = 1957.5
$ws.Range("M122").Value = 492.5

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7840.0977
$ws.Range("I32").Value = 5984.973
$ws.Range("J32").Value = 25000
$ws.Range("K32").Value = 5984.973
$ws.Range("L32").Value = 25000
$ws.Range("M32").Value = -5697.973
$ws.Range("N32").Value = -25574
$ws.Range("H45").Value = 1696.4
$ws.Range("I45").Value = 1696.4
$ws.Range("K45").Value = 1696.4
$ws.Range("M45").Value = -1319.4
$ws.Range("H63").Value = 5792.2856
$ws.Range("J63").Value = 7661.625
$ws.Range("L63").Value = 7661.625
$ws.Range("N63").Value = -9033.625
$ws.Range("H66").Value = 5792.2856
$ws.Range("J66").Value = 7661.625
$ws.Range("L66").Value = 38308.125
$ws.Range("N66").Value = -45172.125
$ws.Range("H74").Value = 1741.4259
$ws.Range("I74").Value = 1150.1666
$ws.Range("K74").Value = 1150.1666
$ws.Range("M74").Value = -276.1666
$ws.Range("H77").Value = 1741.4259
$ws.Range("I77").Value = 1150.1666
$ws.Range("K77").Value = 5750.833000000001
$ws.Range("M77").Value = -1382.833000000001
$ws.Range("H110").Value = 6599.4
$ws.Range("I110").Value = 6599.4
$ws.Range("K110").Value = 6599.4
$ws.Range("M110").Value = -4554.4
$ws.Range("H120").Value = 58999
$ws.Range("I120").Value = 58999
$ws.Range("K120").Value = 58999
$ws.Range("M120").Value = -54161
$ws.Range("H122").Value = 2452.44
$ws.Range("I122").Value = 1638.1875
$ws.Range("J122").Value = 3900
$ws.Range("K122").Value = 4914.5625
$ws.Range("L122").Value = 11700
$ws.Range("M122").Value = -2464.5625
$ws.Range("N122").Value = -16600
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").Value = ""

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 963.84
$ws.Range("I86").Value = 917.36365
$ws.Range("K86").Value = 917.36365
$ws.Range("M86").Value = 205.63635
$ws.Range("H89").Value = 963.84
$ws.Range("I89").Value = 917.36365
$ws.Range("K89").Value = 4586.81825
$ws.Range("M89").Value = 1029.18175
$ws.Range("H99").Value = 3835.7144
$ws.Range("I99").Value = 3792.9092
$ws.Range("J99").Value = 3992.6667
$ws.Range("K99").Value = 3792.9092
$ws.Range("L99").Value = 3992.6667
$ws.Range("M99").Value = -2294.9092
$ws.Range("N99").Value = -6988.6667
$ws.Range("H105").Value = 2961.125
$ws.Range("I105").Value = 2128.5667
$ws.Range("K105").Value = 2128.5667
$ws.Range("M105").Value = -381.5666999999999
$ws.Range("H134").Value = 4095.4167
$ws.Range("I134").Value = 2760
$ws.Range("K134").Value = 8280
$ws.Range("M134").Value = -5745

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4176.5
$ws.Range("I31").Value = 2624.9048
$ws.Range("J31").Value = 7138.636
$ws.Range("K31").Value = 2624.9048
$ws.Range("L31").Value = 7138.636
$ws.Range("M31").Value = -2329.9048
$ws.Range("N31").Value = -7728.636
$ws.Range("H34").Value = 4176.5
$ws.Range("I34").Value = 2624.9048
$ws.Range("J34").Value = 7138.636
$ws.Range("K34").Value = 2624.9048
$ws.Range("L34").Value = 7138.636
$ws.Range("M34").Value = -2422.9048
$ws.Range("N34").Value = -7542.636
$ws.Range("H105").Value = 446.4
$ws.Range("I105").Value = 446.4
$ws.Range("K105").Value = 446.4
$ws.Range("M105").Value = 1300.6
$ws.Range("H122").Value = 3098.7693
$ws.Range("I122").Value = 2898.818
$ws.Range("K122").Value = 8696.454000000002
$ws.Range("M122").Value = -6246.454000000002
$ws.Range("H132").Value = 7014
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 7014
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 21042
$ws.Range("M132").Value = ""
$ws.Range("N132").Value = -26102
$ws.Range("H134").Value = 2637.75
$ws.Range("I134").Value = 1946.5714
$ws.Range("J134").Value = 4711.2856
$ws.Range("K134").Value = 5839.7142
$ws.Range("L134").Value = 14133.8568
$ws.Range("M134").Value = -3304.7142
$ws.Range("N134").Value = -19203.8568

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 2415305.8
$ws.Range("J4").Value = 2700
$ws.Range("L4").Value = 8100
$ws.Range("N4").Value = -8324
$ws.Range("H98").Value = 705.8333
$ws.Range("I98").Value = 697.5
$ws.Range("J98").Value = 722.5
$ws.Range("K98").Value = 2092.5
$ws.Range("L98").Value = 2167.5
$ws.Range("M98").Value = -594.5
$ws.Range("N98").Value = -5163.5
$ws.Range("H122").Value = 311.42856
$ws.Range("I122").Value = 277.3
$ws.Range("J122").Value = 396.75
$ws.Range("K122").Value = 2495.7
$ws.Range("L122").Value = 3570.75
$ws.Range("M122").Value = -45.70000000000027
$ws.Range("N122").Value = -8470.75

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4750.8887
$ws.Range("I70").Value = 4504
$ws.Range("K70").Value = 4504
$ws.Range("M70").Value = -4234
$ws.Range("H73").Value = 4750.8887
$ws.Range("I73").Value = 4504
$ws.Range("K73").Value = 4504
$ws.Range("M73").Value = -3568
$ws.Range("H122").Value = 114500.336
$ws.Range("I122").Value = 2551
$ws.Range("K122").Value = 7653
$ws.Range("M122").Value = -5203

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H13").Value = 405
$ws.Range("I13").Value = 405
$ws.Range("K13").Value = 405
$ws.Range("M13").Value = -265
$ws.Range("H40").Value = 1858.4
$ws.Range("I40").Value = 1848
$ws.Range("K40").Value = 1848
$ws.Range("M40").Value = -1712
$ws.Range("H82").Value = 3234.9473
$ws.Range("I82").Value = 3364.8333
$ws.Range("K82").Value = 3364.8333
$ws.Range("M82").Value = -3003.8333
$ws.Range("H85").Value = 3234.9473
$ws.Range("I85").Value = 3364.8333
$ws.Range("K85").Value = 3364.8333
$ws.Range("M85").Value = -2116.8333
$ws.Range("H93").Value = 3714.8572
$ws.Range("I93").Value = 3714.8572
$ws.Range("K93").Value = 3714.8572
$ws.Range("M93").Value = -2466.8572
$ws.Range("H100").Value = 2250
$ws.Range("J100").Value = 3000
$ws.Range("L100").Value = 3000
$ws.Range("N100").Value = -4082
$ws.Range("H122").Value = 4039.4546
$ws.Range("I122").Value = 4043.5
$ws.Range("K122").Value = 12130.5
$ws.Range("M122").Value = -9680.5
$ws.Range("H133").Value = 33333
$ws.Range("J133").Value = 33333
$ws.Range("L133").Value = 33333
$ws.Range("N133").Value = -38393
$ws.Range("H136").Value = 3199.5
$ws.Range("I136").Value = 3199.5
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 9598.5
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -7048.5
$ws.Range("N136").Value = ""

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2986
$ws.Range("I122").Value = 3490.5386
$ws.Range("K122").Value = 10471.6158
$ws.Range("M122").Value = -8021.6158
